$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 295: D295 and F295 values changed ---
$ws.Cells.Item(295, 4).Value = 36.891   # D295 (high)
$ws.Cells.Item(295, 6).Value = 36.57    # F295 (close)

# --- Append new rows 296-298, cloning the date-cell format from A295 ---
$newRows = @(
    @{ Row = 296; DateSerial = 45170.33333333334; Open = 36.57; High = 36.57; Low = 36.57; Close = 36.57; Volume = 0 },
    @{ Row = 297; DateSerial = 45201.375;          Open = 36.57; High = 36.57; Low = 36.57; Close = 36.57; Volume = 0 },
    @{ Row = 298; DateSerial = 45231.375;          Open = 36.57; High = 36.57; Low = 36.57; Close = 36.57; Volume = 0 }
)

foreach ($r in $newRows) {
    $rowIdx = $r.Row

    # Clone formatting (number format / font / border / alignment) from A295
    $ws.Cells.Item(295, 1).Copy($ws.Cells.Item($rowIdx, 1))
    $ws.Cells.Item($rowIdx, 1).Value = $r.DateSerial

    $ws.Cells.Item($rowIdx, 2).Value = "FX_IDC:USDUAH"
    $ws.Cells.Item($rowIdx, 3).Value = $r.Open
    $ws.Cells.Item($rowIdx, 4).Value = $r.High
    $ws.Cells.Item($rowIdx, 5).Value = $r.Low
    $ws.Cells.Item($rowIdx, 6).Value = $r.Close
    $ws.Cells.Item($rowIdx, 7).Value = $r.Volume
}
